$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.104986190795898
$ws.Range("B1").Value = 1.911372423171997
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.178259372711182
$ws.Range("E1").Value = 1.222580432891846
